# Auto-generated edit script: updates market-price derived columns
# (currentAveragePrice* / LevePrice* / LeveProfit*) for specific leve rows
# across several job sheets, per scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 35000
$ws.Range("J13").Value = 35000
$ws.Range("L13").Value = 35000
$ws.Range("N13").Value = -35338

$ws.Range("H33").Value = 161.4
$ws.Range("I33").Value = 170.61539
$ws.Range("J33").Value = 101.5
$ws.Range("K33").Value = 170.61539
$ws.Range("L33").Value = 101.5
$ws.Range("M33").Value = 58.38461000000001
$ws.Range("N33").Value = -559.5

$ws.Range("H113").Value = 2511.2307
$ws.Range("I113").Value = 2525.5789
$ws.Range("J113").Value = 2472.2856
$ws.Range("K113").Value = 2525.5789
$ws.Range("L113").Value = 2472.2856
$ws.Range("M113").Value = 728.4211
$ws.Range("N113").Value = -8980.285599999999

$ws.Range("H116").Value = 4299.923
$ws.Range("I116").Value = 5559.8
$ws.Range("J116").Value = 3512.5
$ws.Range("K116").Value = 5559.8
$ws.Range("L116").Value = 3512.5
$ws.Range("M116").Value = -2117.8
$ws.Range("N116").Value = -10396.5

$ws.Range("H132").Value = 1792.2858
$ws.Range("I132").Value = 1963.7
$ws.Range("J132").Value = 763.8
$ws.Range("K132").Value = 5891.1
$ws.Range("L132").Value = 2291.4
$ws.Range("M132").Value = -3361.1
$ws.Range("N132").Value = -7351.4

$ws.Range("H137").Value = 763.3469
$ws.Range("I137").Value = 709.2632
$ws.Range("K137").Value = 2127.7896
$ws.Range("M137").Value = 422.2103999999999

$ws.Range("H138").Value = 910.35
$ws.Range("I138").Value = 533.27026
$ws.Range("J138").Value = 1983.5769
$ws.Range("K138").Value = 1599.81078
$ws.Range("L138").Value = 5950.7307
$ws.Range("M138").Value = 3540.18922
$ws.Range("N138").Value = -16230.7307

$ws.Range("H141").Value = 1884.7021
$ws.Range("I141").Value = 628.3158
$ws.Range("J141").Value = 7189.4443
$ws.Range("K141").Value = 1884.9474
$ws.Range("L141").Value = 21568.3329
$ws.Range("M141").Value = 3295.0526
$ws.Range("N141").Value = -31928.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 10000
$ws.Range("J8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("N8").Value = -10288

$ws.Range("H61").Value = 1087.8478
$ws.Range("I61").Value = 721.7632
$ws.Range("J61").Value = 2826.75
$ws.Range("K61").Value = 721.7632
$ws.Range("L61").Value = 2826.75
$ws.Range("M61").Value = -509.7632
$ws.Range("N61").Value = -3250.75

$ws.Range("H74").Value = 785.35297
$ws.Range("I74").Value = 698.5965
$ws.Range("J74").Value = 1234.909
$ws.Range("K74").Value = 698.5965
$ws.Range("L74").Value = 1234.909
$ws.Range("M74").Value = 175.4035
$ws.Range("N74").Value = -2982.909

$ws.Range("H77").Value = 785.35297
$ws.Range("I77").Value = 698.5965
$ws.Range("J77").Value = 1234.909
$ws.Range("K77").Value = 3492.9825
$ws.Range("L77").Value = 6174.545
$ws.Range("M77").Value = 875.0174999999999
$ws.Range("N77").Value = -14910.545

$ws.Range("H136").Value = 1087.8478
$ws.Range("I136").Value = 721.7632
$ws.Range("J136").Value = 2826.75
$ws.Range("K136").Value = 2165.2896
$ws.Range("L136").Value = 8480.25
$ws.Range("M136").Value = 384.7103999999999
$ws.Range("N136").Value = -13580.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 70009
$ws.Range("J14").Value = 70009
$ws.Range("L14").Value = 70009
$ws.Range("N14").Value = -70353

$ws.Range("H86").Value = 2240.12
$ws.Range("I86").Value = 2024.0952
$ws.Range("J86").Value = 3374.25
$ws.Range("K86").Value = 2024.0952
$ws.Range("L86").Value = 3374.25
$ws.Range("M86").Value = -901.0952
$ws.Range("N86").Value = -5620.25

$ws.Range("H89").Value = 2240.12
$ws.Range("I89").Value = 2024.0952
$ws.Range("J89").Value = 3374.25
$ws.Range("K89").Value = 10120.476
$ws.Range("L89").Value = 16871.25
$ws.Range("M89").Value = -4504.476000000001
$ws.Range("N89").Value = -28103.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3909.2307
$ws.Range("I58").Value = 1098.6666
$ws.Range("J58").Value = 13277.777
$ws.Range("K58").Value = 1098.6666
$ws.Range("L58").Value = 13277.777
$ws.Range("M58").Value = -895.6666
$ws.Range("N58").Value = -13683.777

$ws.Range("H99").Value = 2908.25
$ws.Range("I99").Value = 2530
$ws.Range("J99").Value = 3178.4285
$ws.Range("K99").Value = 2530
$ws.Range("L99").Value = 3178.4285
$ws.Range("M99").Value = -1032
$ws.Range("N99").Value = -6174.4285

$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -7900

$ws.Range("H126").Value = 2908.25
$ws.Range("I126").Value = 2530
$ws.Range("J126").Value = 3178.4285
$ws.Range("K126").Value = 7590
$ws.Range("L126").Value = 9535.2855
$ws.Range("M126").Value = -5120
$ws.Range("N126").Value = -14475.2855

$ws.Range("H132").Value = 1963.6061
$ws.Range("I132").Value = 1314.2858
$ws.Range("K132").Value = 3942.8574
$ws.Range("M132").Value = -1412.8574

$ws.Range("H134").Value = 1499.8163
$ws.Range("I134").Value = 1339.45
$ws.Range("J134").Value = 2212.5557
$ws.Range("K134").Value = 4018.35
$ws.Range("L134").Value = 6637.6671
$ws.Range("M134").Value = -1483.35
$ws.Range("N134").Value = -11707.6671

$ws.Range("H136").Value = 3909.2307
$ws.Range("I136").Value = 1098.6666
$ws.Range("J136").Value = 13277.777
$ws.Range("K136").Value = 3295.9998
$ws.Range("L136").Value = 39833.331
$ws.Range("M136").Value = -745.9998000000001
$ws.Range("N136").Value = -44933.331

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1628.0186
$ws.Range("I132").Value = 1426.7333
$ws.Range("J132").Value = 2634.4443
$ws.Range("K132").Value = 4280.199900000001
$ws.Range("L132").Value = 7903.3329
$ws.Range("M132").Value = -1750.199900000001
$ws.Range("N132").Value = -12963.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7859.4375
$ws.Range("I16").Value = 11207.2
$ws.Range("J16").Value = 2279.8333
$ws.Range("K16").Value = 11207.2
$ws.Range("L16").Value = 2279.8333
$ws.Range("M16").Value = -11037.2
$ws.Range("N16").Value = -2619.8333

$ws.Range("H132").Value = 1509.9839
$ws.Range("I132").Value = 1337.9811
$ws.Range("J132").Value = 2522.889
$ws.Range("K132").Value = 4013.9433
$ws.Range("L132").Value = 7568.667
$ws.Range("M132").Value = -1483.9433
$ws.Range("N132").Value = -12628.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 63336.668
$ws.Range("J19").Value = 63336.668
$ws.Range("L19").Value = 63336.668
$ws.Range("N19").Value = -63684.668

$ws.Range("H132").Value = 1570.9166
$ws.Range("I132").Value = 1239
$ws.Range("J132").Value = 2566.6667
$ws.Range("K132").Value = 3717
$ws.Range("L132").Value = 7700.000100000001
$ws.Range("M132").Value = -1187
$ws.Range("N132").Value = -12760.0001

$ws.Range("H136").Value = 618.6923
$ws.Range("I136").Value = 275.90625
$ws.Range("J136").Value = 2185.7144
$ws.Range("K136").Value = 827.71875
$ws.Range("L136").Value = 6557.1432
$ws.Range("M136").Value = 1722.28125
$ws.Range("N136").Value = -11657.1432
